# Update the "Förändrad" (changed) date column (C) for rows 2-19
# from serial date 45224 (2023-10-25) to serial date 45233 (2023-11-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45224) {
        $cell.Value2 = 45233
    }
}
